# Update cryptocurrency price/volume data as scraped on Sun Apr  7 08:20:39 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.395.03'
$ws.Range("E2").Value = '  +2.10%  '

# Row 3
$ws.Range("D3").Value = '3.389.80'
$ws.Range("E3").Value = '  +1.83%  '

# Row 4
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$cell = $ws.Range("D5")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '588.01'
$cell.Style = $savedStyle
$ws.Range("E5").Value = '  +0.85%  '

# Row 6
$cell = $ws.Range("D6")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '180.50'
$cell.Style = $savedStyle
$ws.Range("E6").Value = '  +1.97%  '

# Row 7
$cell = $ws.Range("D7")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = $savedStyle
$ws.Range("E7").Value = '  -0.04%  '

# Row 8
$ws.Range("E8").Value = '  +1.06%  '

# Row 9
$ws.Range("E9").Value = '  +5.79%  '

# Row 10
$cell = $ws.Range("D10")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.591'
$cell.Style = $savedStyle
$ws.Range("E10").Value = '  +1.52%  '

# Row 11
$cell = $ws.Range("D11")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '48.50'
$cell.Style = $savedStyle
$ws.Range("E11").Value = '  +2.46%  '

# Row 12
$cell = $ws.Range("D12")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0000282'
$cell.Style = $savedStyle
$ws.Range("E12").Value = '  +2.94%  '

# Row 13
$cell = $ws.Range("D13")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '677.98'
$cell.Style = $savedStyle
$ws.Range("E13").Value = '  -2.73%  '

# Row 14
$cell = $ws.Range("D14")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '8.62'
$cell.Style = $savedStyle
$ws.Range("E14").Value = '  +2.24%  '

# Row 15
$ws.Range("D15").Value = '3.935.87'
$ws.Range("E15").Value = '  +1.72%  '

# Row 16
$ws.Range("D16").Value = '69.462.57'
$ws.Range("E16").Value = '  +2.15%  '

# Row 17
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$cell = $ws.Range("D17")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.120'
$cell.Style = $savedStyle
$ws.Range("E17").Value = '  +1.86%  '

# Row 18
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.371.37'
$ws.Range("E18").Value = '  +1.02%  '

# Row 19
$cell = $ws.Range("D19")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '17.63'
$cell.Style = $savedStyle
$ws.Range("E19").Value = '  +0.78%  '

# Row 20
$cell = $ws.Range("D20")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '11.28'
$cell.Style = $savedStyle
$ws.Range("E20").Value = '  +2.01%  '

# Row 21
$cell = $ws.Range("D21")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.904'
$cell.Style = $savedStyle
$ws.Range("E21").Value = '  +0.74%  '

# Row 22
$cell = $ws.Range("D22")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.43'
$cell.Style = $savedStyle
$ws.Range("E22").Value = '  +0.64%  '

# Row 23
$cell = $ws.Range("D23")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '17.20'
$cell.Style = $savedStyle
$ws.Range("E23").Value = '  +0.79%  '

# Row 24
$cell = $ws.Range("D24")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '103.46'
$cell.Style = $savedStyle
$ws.Range("E24").Value = '  +3.80%  '

# Row 25
$cell = $ws.Range("D25")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.94'
$cell.Style = $savedStyle
$ws.Range("E25").Value = '  +0.40%  '

# Row 26
$ws.Range("E26").Value = '  +1.16%  '

# Row 27
$cell = $ws.Range("D27")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '9.67'
$cell.Style = $savedStyle
$ws.Range("E27").Value = '  +0.65%  '

# Row 28
$cell = $ws.Range("D28")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '34.11'
$cell.Style = $savedStyle
$ws.Range("E28").Value = '  +2.62%  '

# Row 29
$cell = $ws.Range("D29")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '8.70'
$cell.Style = $savedStyle
$ws.Range("E29").Value = '  +1.45%  '

# Row 30
$cell = $ws.Range("D30")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.03'
$cell.Style = $savedStyle
$ws.Range("E30").Value = '  -0.94%  '

# Row 31
$ws.Range("E31").Value = '  +0.99%  '

# Row 32
$cell = $ws.Range("D32")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '555.82'
$cell.Style = $savedStyle
$ws.Range("E32").Value = '  -2.16%  '

# Row 33
$cell = $ws.Range("D33")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.60'
$cell.Style = $savedStyle
$ws.Range("E33").Value = '  +5.66%  '

# Row 34
$ws.Range("E34").Value = '  +0.67%  '

# Row 35
$cell = $ws.Range("D35")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '58.30'
$cell.Style = $savedStyle
$ws.Range("E35").Value = '  +1.50%  '

# Row 36
$ws.Range("E36").Value = '  +0.11%  '

# Row 37
$ws.Range("D37").Value = '3.687.10'
$ws.Range("E37").Value = '  -0.07%  '

# Row 38
$ws.Range("E38").Value = '  +4.81%  '

# Row 39
$cell = $ws.Range("D39")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '35.30'
$cell.Style = $savedStyle
$ws.Range("E39").Value = '  +1.86%  '

# Row 40
$cell = $ws.Range("D40")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.27'
$cell.Style = $savedStyle
$ws.Range("E40").Value = '  +2.76%  '

# Row 41
$cell = $ws.Range("D41")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.72'
$cell.Style = $savedStyle
$ws.Range("E41").Value = '  +1.96%  '

# Row 42
$ws.Range("D42").Value = '0.0₃0697'
$ws.Range("E42").Value = '  +3.21%  '

# Row 43
$cell = $ws.Range("D43")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.339'
$cell.Style = $savedStyle
$ws.Range("E43").Value = '  +0.75%  '

# Row 44
$cell = $ws.Range("D44")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0425'
$cell.Style = $savedStyle
$ws.Range("E44").Value = '  +4.32%  '

# Row 45
$cell = $ws.Range("D45")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.28'
$cell.Style = $savedStyle
$ws.Range("E45").Value = '  -2.07%  '

# Row 46
$ws.Range("E46").Value = '  -0.64%  '

# Row 47
$ws.Range("E47").Value = '  +0.74%  '

# Row 48
$cell = $ws.Range("D48")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.41'
$cell.Style = $savedStyle
$ws.Range("E48").Value = '  +5.65%  '

# Row 49
$ws.Range("E49").Value = '  +0.00%  '

# Row 50
$cell = $ws.Range("D50")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '133.27'
$cell.Style = $savedStyle
$ws.Range("E50").Value = '  +2.18%  '

# Row 51
$cell = $ws.Range("D51")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.62'
$cell.Style = $savedStyle
$ws.Range("E51").Value = '  +3.91%  '

